# Populate Sheet1 (user list) and Sheet2 (a lookup/reference row) following
# the order that produces the shared-string table seen in the target file.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Row 2: Sourabh Awasthi -------------------------------------------------
$ws1.Range("D2").Value = "sourabh.awasthi@capgemini.com"
$ws1.Range("J2").Value = "pLgknCtnPL]K"

# --- Row 3: Sandipan Deb ----------------------------------------------------
$ws1.Range("D3").Value = "sandipan.deb@capgemini.com"
$ws1.Range("J3").Value = "6xA2pWyl@`$g?"

# Sheet2 reference row, written at this point in the edit history so that
# "Capgemini" lands right after Sandipan's password in the shared strings.
$ws2.Range("C2").Value = "sandipan.deb@capgemini.com"
$ws2.Range("I2").Value = "Capgemini"

# --- Row 4: Biswaji Deb ------------------------------------------------------
$ws1.Range("D4").Value = "biswaji.deb@capgemini.com"
$ws1.Range("J4").Value = "c2h]RMmKmhFz"

# --- Row 5: Debanjan Das ------------------------------------------------------
$ws1.Range("D5").Value = "debanjan.das@capgemini.com"
$ws1.Range("J5").Value = "VN5cEXfU*X)B"

# --- Row 6: Dhiraj Kajari -----------------------------------------------------
$ws1.Range("D6").Value = "dhiraj.kajari@capgemini.com"
$ws1.Range("J6").Value = "tR5i#!sWVo(A"

# --- Row 7: Manoj Kumar B S ---------------------------------------------------
$ws1.Range("D7").Value = "manoj-kumar.b.s@capgemini.com"
$ws1.Range("J7").Value = "gvfzVJ+VT?&v"

# --- Row 8: Mayur Bhorkar ------------------------------------------------------
$ws1.Range("D8").Value = "mayur.bhorkar@capgemini.com"
$ws1.Range("J8").Value = "fLUi]A?uP049"

# --- Numeric / boolean helper columns for rows 2-8 -----------------------
$ws1.Range("K2:K8").Value = 80
$ws1.Range("M2:M8").Value = $true

# --- Formulas for rows 2-7 (written as one shared-formula block each) -----
$ws1.Range("A2:A7").Formula = "=PROPER(IFERROR(LEFT(C2,FIND(CHAR(46),C2)-1),C2))"
$ws1.Range("B2:B7").Formula = '=IFERROR(PROPER(RIGHT(C2,LEN(C2)-FIND("@",SUBSTITUTE(C2,".","@",((LEN(C2)-LEN(SUBSTITUTE(C2,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C2:C7").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D2,FIND(CHAR(64),D2)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws1.Range("E2:E7").Formula = "=LEFT(H2,FIND(CHAR(46),H2)-1)"
$ws1.Range("F2:F7").Formula = '=CONCATENATE("ITPartner\",I2)'
$ws1.Range("H2:H7").Formula = "=RIGHT(D2,LEN(D2)-FIND(CHAR(64),D2))"
$ws1.Range("I2:I7").Formula = "=PROPER(E2)"
$ws1.Range("P2:P7").Formula = "=COUNTIF(D:D,D2)"

# --- Formulas for row 8 (kept as individual, non-shared formulas) ---------
$ws1.Range("A8").Formula = "=PROPER(IFERROR(LEFT(C8,FIND(CHAR(46),C8)-1),C8))"
$ws1.Range("B8").Formula = '=IFERROR(PROPER(RIGHT(C8,LEN(C8)-FIND("@",SUBSTITUTE(C8,".","@",((LEN(C8)-LEN(SUBSTITUTE(C8,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C8").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D8,FIND(CHAR(64),D8)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws1.Range("E8").Formula = "=LEFT(H8,FIND(CHAR(46),H8)-1)"
$ws1.Range("F8").Formula = '=CONCATENATE("ITPartner\",I8)'
$ws1.Range("H8").Formula = "=RIGHT(D8,LEN(D8)-FIND(CHAR(64),D8))"
$ws1.Range("I8").Formula = "=PROPER(E8)"
$ws1.Range("P8").Formula = "=COUNTIF(D:D,D8)"
